# Update the "Elapsed Duration(Hrs)" (column G) values on several sheets
# of Active_Outages.xlsx to reflect newly-recalculated outage durations.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "R1"; Row = 2; Value = "3883:34:46" },
    @{ Sheet = "R1"; Row = 3; Value = "23:07:24" },
    @{ Sheet = "R1"; Row = 4; Value = "78679:13:57" },
    @{ Sheet = "R2"; Row = 2; Value = "12064:58:27" },
    @{ Sheet = "R2"; Row = 3; Value = "3194:41:56" },
    @{ Sheet = "R2"; Row = 4; Value = "432:53:30" },
    @{ Sheet = "R4"; Row = 2; Value = "2910:48:16" },
    @{ Sheet = "R4"; Row = 3; Value = "138:00:31" },
    @{ Sheet = "R5"; Row = 2; Value = "384:47:15" },
    @{ Sheet = "R6"; Row = 2; Value = "25:19:33" }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Cells.Item($u.Row, 7).Value = $u.Value
}
